# Update "paises.xlsx" (Pais sheet) with refreshed COVID figures.
# This mirrors a newer data pull: most countries just get updated counts,
# but a few change rank versus their neighbours (Uzbekistan overtakes
# Senegal; Mauritania and Angola each move up one rank), which is why
# some rows below them cascade down by one position. The timestamp
# caption in A1 is also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp caption ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 21:10"

# --- Helper: write a full data row (Pais, Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Rows whose ranking is unchanged, just refreshed totals.
Set-Row 4  @("Estados Unidos", 1806862, 13332, 523064, 1178730, 0, 526, 105068)
Set-Row 10 @("Francia",        188625,  1828,  68268,  91586,   0, 57,  28771)
Set-Row 12 @("India",          181827,  8336,  86936,  89706,   0, 205, 5185)

# Uzbekistan overtakes Senegal in total cases -> the two rows swap places.
Set-Row 78 @("Uzbekistan", 3546, 78,  2783, 749,  0, 0, 14)
Set-Row 79 @("Senegal",    3535, 106, 1761, 1732, 0, 1, 42)

# Mauritania moves up to rank 142 (row 138), pushing Reunion ... Togo
# down by one row each.
Set-Row 138 @("Mauritania",             483, 60, 21,  442, 0, 0, 20)
Set-Row 139 @("Reunion",                471, 1,  411, 59,  0, 0, 1)
Set-Row 140 @("Santo Tome y Principe",  463, 0,  68,  383, 0, 0, 12)
Set-Row 141 @("Guayana Francesa",       450, 0,  172, 277, 0, 0, 1)
Set-Row 142 @("Estado de Palestina",    447, 1,  368, 76,  0, 0, 3)
Set-Row 143 @("Taiwan",                 442, 0,  421, 14,  0, 0, 7)
Set-Row 144 @("Togo",                   428, 0,  202, 213, 0, 0, 13)

# Angola moves up to rank 183 (row 178), pushing Liechtenstein down one row.
Set-Row 178 @("Angola",       84, 3, 18, 62, 0, 0, 4)
Set-Row 179 @("Liechtenstein",82, 0, 55, 26, 0, 0, 1)
